# Fix a typo in the "Unkown" row's placeholder token and populate the
# previously-empty placeholder for the "Electric" living-expense row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# H5: {unkown} -> {unknown}
$ws.Cells.Item(5, 8).Value = "{unknown}"

# B13: (empty) -> {electricity}
$ws.Cells.Item(13, 2).Value = "{electricity}"

Write-Output "applied placeholder fixes"
